# Error Calculations and Plots
# Apply the edits described in the commit:
#  - Remove the "RM 232" row and the "SC 92" row from the data table
#    (rows shift up to close the gaps).
#  - Fill in / clear a handful of previously-missing / now-missing
#    values in column E (and a couple in column D) as re-measured.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the two rows that were dropped from the dataset.
#    Row 26 = "RM 232". After it is removed, the row that used to be 28
#    ("SC 92") becomes row 27, so delete that next.
$ws.Range("A26:F26").EntireRow.Delete()
$ws.Range("A27:F27").EntireRow.Delete()

# 2) Apply the individual cell corrections (final row numbers, after the
#    two rows above were removed).
$ws.Range("E6").Value = -5.7
$ws.Range("E8").ClearContents()
$ws.Range("E19").Value = -6.5
$ws.Range("E21").ClearContents()
$ws.Range("E23").Value = -7

$ws.Range("D26").ClearContents()
$ws.Range("D27").Value = -14.6
$ws.Range("E27").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("E29").Value = -6.8
